$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "gemini-1.5-pro"
$ws.Range("C9").Value = "-0.0 ± 0.61"
$ws.Range("D9").Value = "0.46 ± 0.76"
$ws.Range("J9").Value = "0.82 ± 0.15"
$ws.Range("K9").Value = "0.84 ± 0.15"
$ws.Range("L9").Value = "0.83 ± 0.15"
$ws.Range("M9").Value = "0.89 ± 0.17"
$ws.Range("N9").Value = "0.97 ± 0.17"
$ws.Range("P9").Value = "0.5 ± 0.16"
$ws.Range("Q9").Value = "2.75 ± 1.54"
$ws.Range("R9").Value = "0.032 ± 0.00"
$ws.Range("S9").Value = "0.92 ± 0.17"
$ws.Range("T9").Value = "0.93 ± 0.21"
$ws.Range("U9").Value = "2.78 ± 1.17"
$ws.Range("V9").Value = "0.63 ± 0.4"
$ws.Range("W9").Value = "0.92 ± 0.17"
$ws.Range("X9").Value = "1.24 ± 0.31"

# Force-create placeholder cells for the columns that are blank in this
# row (matching the source data which keeps an empty cell record for
# every column instead of leaving a true gap).
$ws.Range("B9").Style = "Normal"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").Style = "Normal"
$ws.Range("I9").Style = "Normal"
$ws.Range("O9").Style = "Normal"
$ws.Range("Y9").Style = "Normal"
